# Matriz de trazabilidad - "Modificar Usuario Repartidor + Actualizacion Matriz de Trazabilidad"
#
# 1) Row 33 (item #31) changes its "Estado" (H33) from "Falta" to "Completado",
#    which also swaps the red "Falta" fill/font for the green "Completado"
#    fill/font already used elsewhere in the sheet (e.g. H6/H34).
# 2) The "Fecha de estado" (column I) advances one day, from 2020-10-31
#    (serial 44135) to 2020-11-01 (serial 44136), for rows 5, 7-15, 33, 35-39.
# 3) The active selection moves from H34 to I5 (and the stale scrolled-down
#    view snaps back in the process).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) H33: "Falta" -> "Completado", reusing the existing "Completado" look ---
$ws.Range("H33").Value = "Completado"
$ws.Range("H6").Copy()
$ws.Range("H33").PasteSpecial(-4122)   # xlPasteFormats: copy only fill/font/etc.
$excel.CutCopyMode = $false

# --- 2) Bump "Fecha de estado" from 44135 to 44136 for the affected rows ---
$rows = @(5, 7, 8, 9, 10, 11, 12, 13, 14, 15, 33, 35, 36, 37, 38, 39)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 9).Value = 44136   # column I = 9
}

# --- 3) Move the selection/active cell to I5 ---
$ws.Range("I5").Select()
